# Updates the crypto price (column D) and 1h volume/change (column E) values
# for rows 2-51 of the active worksheet, reflecting the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.789.19"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.754.16"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'602.61"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'169.70"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "3.753.52"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "'6.36"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "'38.31"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "4.381.48"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "3.750.94"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "68.802.37"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "'7.32"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'17.19"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").Value = "'10.79"
$ws.Range("E21").Value = "  +16.20%  "
$ws.Range("D22").Value = "'494.91"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'0.731"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Value = "'85.53"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'0.0000147"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "'12.45"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'2.54"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("D31").Value = "'2.99"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "'7.97"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "'32.11"
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").Value = "3.899.21"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").Value = "3.687.92"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").Value = "'5.86"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.133"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "'0.327"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").Value = "'439.52"
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("D43").Value = "'49.00"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'8.52"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'40.79"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "2.822.27"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "'141.20"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").Value = "'0.0356"
$ws.Range("E51").Value = "  +1.31%  "
